# 自动更新价格数据: insert a new latest-date row at the top of the data
# table (row 2, just below the header), pushing all existing rows down
# by one. The new row carries the same constant metric values as every
# other row (783.5 / 1112 / 3610) with the new date 2026-02-18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 2..90 down to 3..91, creating a new blank row 2.
$ws.Rows.Item(2).Insert()

# Column A holds dates stored as plain text (e.g. "2026-02-17"). Excel
# would otherwise auto-convert a "YYYY-MM-DD"-looking literal into a
# date serial, so force the cell to Text before assigning, then drop
# back to the workbook's Normal style so no stray formatting lingers.
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "2026-02-18"
$ws.Cells.Item(2, 1).Style = "Normal"

$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610
